$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Cells.Item(32, 8).Value = 21855.438  # H32: 25691.928 -> 21855.438
$ws.Cells.Item(32, 10).Value = 19381.637  # J32: 24799.777 -> 19381.637
$ws.Cells.Item(32, 12).Value = 19381.637  # L32: 24799.777 -> 19381.637
$ws.Cells.Item(32, 14).Value = -20033.637  # N32: -25451.777 -> -20033.637

# Row 39
$ws.Cells.Item(39, 8).Value = 257.75  # H39: 460.625 -> 257.75
$ws.Cells.Item(39, 9).Value = 60.333332  # I39: 74.25 -> 60.333332
$ws.Cells.Item(39, 10).Value = 850  # J39: 847 -> 850
$ws.Cells.Item(39, 11).Value = 180.999996  # K39: 222.75 -> 180.999996
$ws.Cells.Item(39, 12).Value = 2550  # L39: 2541 -> 2550
$ws.Cells.Item(39, 13).Value = 115.000004  # M39: 73.25 -> 115.000004
$ws.Cells.Item(39, 14).Value = -3142  # N39: -3133 -> -3142

# Row 58
$ws.Cells.Item(58, 8).Value = 378.23077  # H58: 335.07693 -> 378.23077
$ws.Cells.Item(58, 9).Value = 268.08334  # I58: 303.9 -> 268.08334
$ws.Cells.Item(58, 10).Value = 1700  # J58: 439 -> 1700
$ws.Cells.Item(58, 11).Value = 804.2500200000001  # K58: 911.6999999999999 -> 804.2500200000001
$ws.Cells.Item(58, 12).Value = 5100  # L58: 1317 -> 5100
$ws.Cells.Item(58, 13).Value = -654.2500200000001  # M58: -761.6999999999999 -> -654.2500200000001
$ws.Cells.Item(58, 14).Value = -5400  # N58: -1617 -> -5400

# Row 92
$ws.Cells.Item(92, 8).Value = 11017.8  # H92: 15356.143 -> 11017.8
$ws.Cells.Item(92, 9).Value = 11017.8  # I92: 15356.143 -> 11017.8
$ws.Cells.Item(92, 11).Value = 11017.8  # K92: 15356.143 -> 11017.8
$ws.Cells.Item(92, 13).Value = -9769.799999999999  # M92: -14108.143 -> -9769.799999999999

# Row 137
$ws.Cells.Item(137, 8).Value = 1273.95  # H137: 1377.6904 -> 1273.95
$ws.Cells.Item(137, 9).Value = 1225.8  # I137: 1208.5938 -> 1225.8
$ws.Cells.Item(137, 10).Value = 1418.4  # J137: 1918.8 -> 1418.4
$ws.Cells.Item(137, 11).Value = 3677.4  # K137: 3625.7814 -> 3677.4
$ws.Cells.Item(137, 12).Value = 4255.200000000001  # L137: 5756.4 -> 4255.200000000001
$ws.Cells.Item(137, 13).Value = -1127.4  # M137: -1075.7814 -> -1127.4
$ws.Cells.Item(137, 14).Value = -9355.200000000001  # N137: -10856.4 -> -9355.200000000001

# Row 138
$ws.Cells.Item(138, 8).Value = 3171.7544  # H138: 3184.7856 -> 3171.7544
$ws.Cells.Item(138, 10).Value = 4739.0645  # J138: 4815.6333 -> 4739.0645
$ws.Cells.Item(138, 12).Value = 14217.1935  # L138: 14446.8999 -> 14217.1935
$ws.Cells.Item(138, 14).Value = -24497.1935  # N138: -24726.8999 -> -24497.1935


# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Cells.Item(16, 8).Value = 500  # H16: 544.4 -> 500
$ws.Cells.Item(16, 9).Value = 500  # I16: 544.4 -> 500
$ws.Cells.Item(16, 11).Value = 500  # K16: 544.4 -> 500
$ws.Cells.Item(16, 13).Value = -213  # M16: -257.4 -> -213

# Row 55
$ws.Cells.Item(55, 8).Value = 13203.714  # H55: 10238 -> 13203.714
$ws.Cells.Item(55, 9).Value = 8107.25  # I55: 8285.799999999999 -> 8107.25
$ws.Cells.Item(55, 11).Value = 8107.25  # K55: 8285.799999999999 -> 8107.25
$ws.Cells.Item(55, 13).Value = -7792.25  # M55: -7970.799999999999 -> -7792.25

# Row 74
$ws.Cells.Item(74, 8).Value = 8203.23  # H74: 7935.815 -> 8203.23
$ws.Cells.Item(74, 9).Value = 1158.1052  # I74: 1149.35 -> 1158.1052
$ws.Cells.Item(74, 11).Value = 1158.1052  # K74: 1149.35 -> 1158.1052
$ws.Cells.Item(74, 13).Value = -284.1052  # M74: -275.3499999999999 -> -284.1052

# Row 77
$ws.Cells.Item(77, 8).Value = 8203.23  # H77: 7935.815 -> 8203.23
$ws.Cells.Item(77, 9).Value = 1158.1052  # I77: 1149.35 -> 1158.1052
$ws.Cells.Item(77, 11).Value = 5790.526  # K77: 5746.75 -> 5790.526
$ws.Cells.Item(77, 13).Value = -1422.526  # M77: -1378.75 -> -1422.526

# Row 102
$ws.Cells.Item(102, 8).Value = 3265.16  # H102: 3146.8462 -> 3265.16
$ws.Cells.Item(102, 9).Value = 3164.9546  # I102: 3035.5652 -> 3164.9546
$ws.Cells.Item(102, 11).Value = 3164.9546  # K102: 3035.5652 -> 3164.9546
$ws.Cells.Item(102, 13).Value = -1542.9546  # M102: -1413.5652 -> -1542.9546

# Row 110
$ws.Cells.Item(110, 8).Value = 7938.4165  # H110: 8100.778 -> 7938.4165
$ws.Cells.Item(110, 9).Value = 11687.353  # I110: 11418.333 -> 11687.353
$ws.Cells.Item(110, 10).Value = 4584.1055  # J110: 4783.222 -> 4584.1055
$ws.Cells.Item(110, 11).Value = 11687.353  # K110: 11418.333 -> 11687.353
$ws.Cells.Item(110, 12).Value = 4584.1055  # L110: 4783.222 -> 4584.1055
$ws.Cells.Item(110, 13).Value = -9642.352999999999  # M110: -9373.333000000001 -> -9642.352999999999
$ws.Cells.Item(110, 14).Value = -8674.1055  # N110: -8873.222 -> -8674.1055

# Row 132
$ws.Cells.Item(132, 8).Value = 3647.6843  # H132: 3806.4443 -> 3647.6843
$ws.Cells.Item(132, 9).Value = 1881.2142  # I132: 1965.1538 -> 1881.2142
$ws.Cells.Item(132, 11).Value = 5643.642599999999  # K132: 5895.4614 -> 5643.642599999999
$ws.Cells.Item(132, 13).Value = -3113.642599999999  # M132: -3365.4614 -> -3113.642599999999


# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 44
$ws.Cells.Item(44, 8).Value = 68348.336  # H44: 90015 -> 68348.336
$ws.Cells.Item(44, 10).Value = 62500  # J44: 95000 -> 62500
$ws.Cells.Item(44, 12).Value = 62500  # L44: 95000 -> 62500
$ws.Cells.Item(44, 14).Value = -63494  # N44: -95994 -> -63494

# Row 76
$ws.Cells.Item(76, 8).Value = 19999.666  # H76: 19224.75 -> 19999.666
$ws.Cells.Item(76, 10).Value = 19999.666  # J76: 19224.75 -> 19999.666
$ws.Cells.Item(76, 12).Value = 19999.666  # L76: 19224.75 -> 19999.666
$ws.Cells.Item(76, 14).Value = -20629.666  # N76: -19854.75 -> -20629.666

# Row 79
$ws.Cells.Item(79, 8).Value = 19999.666  # H79: 19224.75 -> 19999.666
$ws.Cells.Item(79, 10).Value = 19999.666  # J79: 19224.75 -> 19999.666
$ws.Cells.Item(79, 12).Value = 19999.666  # L79: 19224.75 -> 19999.666
$ws.Cells.Item(79, 14).Value = -22183.666  # N79: -21408.75 -> -22183.666

# Row 86
$ws.Cells.Item(86, 8).Value = 3849.1428  # H86: 3646.3333 -> 3849.1428
$ws.Cells.Item(86, 9).Value = 2741.2856  # I86: 2386.111 -> 2741.2856
$ws.Cells.Item(86, 10).Value = 4957  # J86: 5536.6665 -> 4957
$ws.Cells.Item(86, 11).Value = 2741.2856  # K86: 2386.111 -> 2741.2856
$ws.Cells.Item(86, 12).Value = 4957  # L86: 5536.6665 -> 4957
$ws.Cells.Item(86, 13).Value = -1618.2856  # M86: -1263.111 -> -1618.2856
$ws.Cells.Item(86, 14).Value = -7203  # N86: -7782.6665 -> -7203

# Row 89
$ws.Cells.Item(89, 8).Value = 3849.1428  # H89: 3646.3333 -> 3849.1428
$ws.Cells.Item(89, 9).Value = 2741.2856  # I89: 2386.111 -> 2741.2856
$ws.Cells.Item(89, 10).Value = 4957  # J89: 5536.6665 -> 4957
$ws.Cells.Item(89, 11).Value = 13706.428  # K89: 11930.555 -> 13706.428
$ws.Cells.Item(89, 12).Value = 24785  # L89: 27683.3325 -> 24785
$ws.Cells.Item(89, 13).Value = -8090.428  # M89: -6314.555 -> -8090.428
$ws.Cells.Item(89, 14).Value = -36017  # N89: -38915.3325 -> -36017

# Row 105
$ws.Cells.Item(105, 8).Value = 4033.3333  # H105: 2878.4 -> 4033.3333
$ws.Cells.Item(105, 9).Value = 6722.25  # I105: 3076.9092 -> 6722.25
$ws.Cells.Item(105, 10).Value = 2688.875  # J105: 2635.7778 -> 2688.875
$ws.Cells.Item(105, 11).Value = 6722.25  # K105: 3076.9092 -> 6722.25
$ws.Cells.Item(105, 12).Value = 2688.875  # L105: 2635.7778 -> 2688.875
$ws.Cells.Item(105, 13).Value = -4975.25  # M105: -1329.9092 -> -4975.25
$ws.Cells.Item(105, 14).Value = -6182.875  # N105: -6129.7778 -> -6182.875

# Row 134
$ws.Cells.Item(134, 8).Value = 1581.8096  # H134: 1692.186 -> 1581.8096
$ws.Cells.Item(134, 9).Value = 1344.5128  # I134: 1455.4872 -> 1344.5128
$ws.Cells.Item(134, 10).Value = 4666.6665  # J134: 4000 -> 4666.6665
$ws.Cells.Item(134, 11).Value = 4033.5384  # K134: 4366.461600000001 -> 4033.5384
$ws.Cells.Item(134, 12).Value = 13999.9995  # L134: 12000 -> 13999.9995
$ws.Cells.Item(134, 13).Value = -1498.5384  # M134: -1831.461600000001 -> -1498.5384
$ws.Cells.Item(134, 14).Value = -19069.9995  # N134: -17070 -> -19069.9995


# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 994.26666  # H16: 872.3158 -> 994.26666
$ws.Cells.Item(16, 9).Value = 817.2222  # I16: 693.46155 -> 817.2222
$ws.Cells.Item(16, 11).Value = 817.2222  # K16: 693.46155 -> 817.2222
$ws.Cells.Item(16, 13).Value = -530.2222  # M16: -406.46155 -> -530.2222

# Row 45
$ws.Cells.Item(45, 8).Value = 11999  # H45: 0 -> 11999
$ws.Cells.Item(45, 9).Value = 11999  # I45: 0 -> 11999
$ws.Cells.Item(45, 11).Value = 11999  # K45: 0 -> 11999
$ws.Cells.Item(45, 13).Value = -11406  # M45: None -> -11406

# Row 60
$ws.Cells.Item(60, 8).Value = 15218.6  # H60: 16523.25 -> 15218.6
$ws.Cells.Item(60, 9).Value = 15218.6  # I60: 16523.25 -> 15218.6
$ws.Cells.Item(60, 11).Value = 15218.6  # K60: 16523.25 -> 15218.6
$ws.Cells.Item(60, 13).Value = -14707.6  # M60: -16012.25 -> -14707.6

# Row 113
$ws.Cells.Item(113, 8).Value = 994.26666  # H113: 872.3158 -> 994.26666
$ws.Cells.Item(113, 9).Value = 817.2222  # I113: 693.46155 -> 817.2222
$ws.Cells.Item(113, 11).Value = 817.2222  # K113: 693.46155 -> 817.2222
$ws.Cells.Item(113, 13).Value = 1352.7778  # M113: 1476.53845 -> 1352.7778

# Row 122
$ws.Cells.Item(122, 8).Value = 918.7059  # H122: 919 -> 918.7059
$ws.Cells.Item(122, 9).Value = 851.125  # I122: 851.4375 -> 851.125
$ws.Cells.Item(122, 11).Value = 2553.375  # K122: 2554.3125 -> 2553.375
$ws.Cells.Item(122, 13).Value = -103.375  # M122: -104.3125 -> -103.375


# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Cells.Item(4, 8).Value = 23535846  # H4: 23084932 -> 23535846
$ws.Cells.Item(4, 9).Value = 26795252  # I4: 26208772 -> 26795252
$ws.Cells.Item(4, 11).Value = 80385756  # K4: 78626316 -> 80385756
$ws.Cells.Item(4, 13).Value = -80385644  # M4: -78626204 -> -80385644

# Row 16
$ws.Cells.Item(16, 8).Value = 350  # H16: 400 -> 350
$ws.Cells.Item(16, 9).Value = 350  # I16: 500 -> 350
$ws.Cells.Item(16, 10).Value = 0  # J16: 350 -> 0
$ws.Cells.Item(16, 11).Value = 1050  # K16: 1500 -> 1050
$ws.Cells.Item(16, 12).Value = 0  # L16: 1050 -> 0
$ws.Cells.Item(16, 13).Value = -877  # M16: -1327 -> -877
$ws.Cells.Item(16, 14).ClearContents()  # N16: -1396 -> (removed)

# Row 17
$ws.Cells.Item(17, 8).Value = 95.40000000000001  # H17: 93 -> 95.40000000000001
$ws.Cells.Item(17, 9).Value = 37.666668  # I17: 33.666668 -> 37.666668
$ws.Cells.Item(17, 11).Value = 113.000004  # K17: 101.000004 -> 113.000004
$ws.Cells.Item(17, 13).Value = 55.999996  # M17: 67.999996 -> 55.999996

# Row 34
$ws.Cells.Item(34, 8).Value = 4275596  # H34: 4632084.5 -> 4275596
$ws.Cells.Item(34, 9).Value = 49.5  # I34: 116.25 -> 49.5
$ws.Cells.Item(34, 10).Value = 6175839  # J34: 6948069 -> 6175839
$ws.Cells.Item(34, 11).Value = 148.5  # K34: 348.75 -> 148.5
$ws.Cells.Item(34, 12).Value = 18527517  # L34: 20844207 -> 18527517
$ws.Cells.Item(34, 13).Value = -64.5  # M34: -264.75 -> -64.5
$ws.Cells.Item(34, 14).Value = -18527685  # N34: -20844375 -> -18527685


# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 47
$ws.Cells.Item(47, 8).Value = 30750  # H47: 30998.75 -> 30750
$ws.Cells.Item(47, 10).Value = 30750  # J47: 30998.75 -> 30750
$ws.Cells.Item(47, 12).Value = 30750  # L47: 30998.75 -> 30750
$ws.Cells.Item(47, 14).Value = -31886  # N47: -32134.75 -> -31886

# Row 70
$ws.Cells.Item(70, 8).Value = 21021.666  # H70: 10859.2 -> 21021.666
$ws.Cells.Item(70, 9).Value = 0  # I70: 7878.2856 -> 0
$ws.Cells.Item(70, 10).Value = 21021.666  # J70: 13467.5 -> 21021.666
$ws.Cells.Item(70, 11).Value = 0  # K70: 7878.2856 -> 0
$ws.Cells.Item(70, 12).Value = 21021.666  # L70: 13467.5 -> 21021.666
$ws.Cells.Item(70, 13).ClearContents()  # M70: -7608.2856 -> (removed)
$ws.Cells.Item(70, 14).Value = -21561.666  # N70: -14007.5 -> -21561.666

# Row 73
$ws.Cells.Item(73, 8).Value = 21021.666  # H73: 10859.2 -> 21021.666
$ws.Cells.Item(73, 9).Value = 0  # I73: 7878.2856 -> 0
$ws.Cells.Item(73, 10).Value = 21021.666  # J73: 13467.5 -> 21021.666
$ws.Cells.Item(73, 11).Value = 0  # K73: 7878.2856 -> 0
$ws.Cells.Item(73, 12).Value = 21021.666  # L73: 13467.5 -> 21021.666
$ws.Cells.Item(73, 13).ClearContents()  # M73: -6942.2856 -> (removed)
$ws.Cells.Item(73, 14).Value = -22893.666  # N73: -15339.5 -> -22893.666

# Row 124
$ws.Cells.Item(124, 8).Value = 29998  # H124: 29998.875 -> 29998
$ws.Cells.Item(124, 10).Value = 29998  # J124: 29998.875 -> 29998
$ws.Cells.Item(124, 12).Value = 29998  # L124: 29998.875 -> 29998
$ws.Cells.Item(124, 14).Value = -39818  # N124: -39818.875 -> -39818


# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 128
$ws.Cells.Item(128, 8).Value = 90000  # H128: 0 -> 90000
$ws.Cells.Item(128, 10).Value = 90000  # J128: 0 -> 90000
$ws.Cells.Item(128, 12).Value = 90000  # L128: 0 -> 90000
$ws.Cells.Item(128, 14).Value = -99960  # N128: None -> -99960


# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Cells.Item(45, 8).Value = 18006.25  # H45: 16415.5 -> 18006.25
$ws.Cells.Item(45, 9).Value = 0  # I45: 13569 -> 0
$ws.Cells.Item(45, 10).Value = 18006.25  # J45: 16984.8 -> 18006.25
$ws.Cells.Item(45, 11).Value = 0  # K45: 13569 -> 0
$ws.Cells.Item(45, 12).Value = 18006.25  # L45: 16984.8 -> 18006.25
$ws.Cells.Item(45, 13).ClearContents()  # M45: -13078 -> (removed)
$ws.Cells.Item(45, 14).Value = -18988.25  # N45: -17966.8 -> -18988.25

# Row 114
$ws.Cells.Item(114, 8).Value = 39990  # H114: 0 -> 39990
$ws.Cells.Item(114, 10).Value = 39990  # J114: 0 -> 39990
$ws.Cells.Item(114, 12).Value = 39990  # L114: 0 -> 39990
$ws.Cells.Item(114, 14).Value = -48668  # N114: None -> -48668

# Row 122
$ws.Cells.Item(122, 8).Value = 1967.5116  # H122: 2026.9025 -> 1967.5116
$ws.Cells.Item(122, 9).Value = 1946.2307  # I122: 1984.2894 -> 1946.2307
$ws.Cells.Item(122, 10).Value = 2175  # J122: 2566.6667 -> 2175
$ws.Cells.Item(122, 11).Value = 5838.6921  # K122: 5952.8682 -> 5838.6921
$ws.Cells.Item(122, 12).Value = 6525  # L122: 7700.000100000001 -> 6525
$ws.Cells.Item(122, 13).Value = -3388.6921  # M122: -3502.8682 -> -3388.6921
$ws.Cells.Item(122, 14).Value = -11425  # N122: -12600.0001 -> -11425

# Row 136
$ws.Cells.Item(136, 8).Value = 3557.8333  # H136: 3790.3635 -> 3557.8333
$ws.Cells.Item(136, 9).Value = 3799.4443  # I136: 4149.375 -> 3799.4443
$ws.Cells.Item(136, 11).Value = 11398.3329  # K136: 12448.125 -> 11398.3329
$ws.Cells.Item(136, 13).Value = -8848.332900000001  # M136: -9898.125 -> -8848.332900000001

# Row 137
$ws.Cells.Item(137, 8).Value = 95365.8  # H137: 95473.5 -> 95365.8
$ws.Cells.Item(137, 10).Value = 95365.8  # J137: 95473.5 -> 95365.8
$ws.Cells.Item(137, 12).Value = 95365.8  # L137: 95473.5 -> 95365.8
$ws.Cells.Item(137, 14).Value = -105565.8  # N137: -105673.5 -> -105565.8

